# Insert a new data row at row 207 (pushing existing rows 207-297 down to 208-298)
# and populate it with the new "Ajo" price record for Macroferia Regional de Talca.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(207).Insert()

$ws.Range("A207").Value = 5
$ws.Range("B207").Value = "Macroferia Regional de Talca"
$ws.Range("C207").Value = "Maule"
$ws.Range("D207").Value = 44704
$ws.Range("E207").Value = 7
$ws.Range("F207").Value = 100112003
$ws.Range("G207").Value = "Ajo"
$ws.Range("H207").Value = "Chino"
$ws.Range("I207").Value = "Primera"
$ws.Range("J207").Value = 300
$ws.Range("K207").Value = 20000
$ws.Range("L207").Value = 20000
$ws.Range("M207").Value = 20000
$ws.Range("N207").Value = "`$/malla 10 kilos"
$ws.Range("O207").Value = "China"
$ws.Range("P207").Value = 2000
$ws.Range("Q207").Value = 10
$ws.Range("R207").Value = "Hortaliza"
